$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 17572.79729737419
$ws.Range("C2").Value = 37492.19165552104
$ws.Range("D2").Value = 73037.81439821956
$ws.Range("E2").Value = 110553.1488376751

# Row 3
$ws.Range("B3").Value = 180531.7646324224
$ws.Range("C3").Value = 335706.1708952549
$ws.Range("D3").Value = 407537.8372402679
$ws.Range("E3").Value = 453519.3175182287

# Row 4
$ws.Range("B4").Value = 19605.67819354744
$ws.Range("C4").Value = 37265.45449342443
$ws.Range("D4").Value = 59525.38698290462
$ws.Range("E4").Value = 77146.58170007616

# Row 6
$ws.Range("B6").Value = 105920.1881086569
$ws.Range("C6").Value = 135008.4585056588
$ws.Range("D6").Value = 124284.6180696946
$ws.Range("E6").Value = 101977.8290434767

# Row 7
$ws.Range("B7").Value = 10958.60998832631
$ws.Range("C7").Value = 22006.77181678597
$ws.Range("D7").Value = 24307.71823809589
$ws.Range("E7").Value = 26511.39833502457

# Row 9
$ws.Range("B9").Value = 811978.2298396495
$ws.Range("C9").Value = 1304953.204775028
$ws.Range("D9").Value = 1770277.589342636
$ws.Range("E9").Value = 2183717.935484926

# Row 12
$ws.Range("B12").Value = 898896.620596268
$ws.Range("C12").Value = 1063696.378554818
$ws.Range("D12").Value = 948319.3777886492
$ws.Range("E12").Value = 725777.4233716547
